$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("B2").Value = "Rakesh11"
$ws.Range("A2").Value = "AutomationCategory99"
$ws.Range("A3").Value = $null
$ws.Range("B3").Value = $null
$ws.Columns.Item(1).ColumnWidth = 25.140625
$ws.Range("A3").Select()
